# "10th - MB for single stock and added new group"
#
# This MarketBeat-rank style report keeps one column per reporting date
# (newest date left-most in row 1) and one row per analyst/firm. Two new
# reporting dates (Jun_26, Jun_27) are being published, so:
#   - the existing date columns (Jun_17/Jun_15/Jun_13/Jun_10, currently in
#     B:E) slide right into E:H, and the two new dates go into the
#     freed-up B:D columns (B=Jun_27, C/D=Jun_26);
#   - every existing analyst row (2-27) keeps its old B:E ratings exactly
#     where they are and simply grows three more "UN" (unchanged) cells
#     into the new F:H columns;
#   - the one highlighted note that lived in the old right-most column
#     (E18, the Jun_10 column) has to follow that column to its new
#     position (H18), and the cell it vacates goes back to a plain "UN";
#   - and a new analyst group is appended as rows 28-29 (Benchmark,
#     Evercore ISI), matching the 4-column (A:D) shape used before this
#     report ever had E:H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift the header row's existing dates right by three columns,
#        and fill the two new dates into the vacated columns. Read with
#        Value2 (Value's getter is not reliable in this host) before any
#        cell gets overwritten.
$oldB1 = $ws.Range("B1").Value2
$oldC1 = $ws.Range("C1").Value2
$oldD1 = $ws.Range("D1").Value2
$oldE1 = $ws.Range("E1").Value2

$ws.Range("H1").Value = $oldE1
$ws.Range("G1").Value = $oldD1
$ws.Range("F1").Value = $oldC1
$ws.Range("E1").Value = $oldB1

$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# --- 2. Existing analyst rows (2-27): B:E stay put, add "UN" into the
#        three newly-opened F:G:H columns for every row.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 6).Value = "UN"   # F
    $ws.Cells.Item($r, 7).Value = "UN"   # G
    $ws.Cells.Item($r, 8).Value = "UN"   # H
}

# --- 3. Row 18's highlighted note moves from the old Jun_10 column (E)
#        to the new Jun_10 column (H); E18 becomes a normal "UN" cell.
$note = $ws.Range("E18").Value2
$ws.Range("H18").Value = $note
$ws.Range("H18").Interior.Color = $ws.Range("E18").Interior.Color

$ws.Range("E18").Value = "UN"
$ws.Range("E18").ClearFormats()

# --- 4. New analyst group appended at the bottom, same 4-column (A:D)
#        shape the sheet used before it grew the E:H date columns.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"

# --- 5. Match the new column widths (C:H all 8.0 characters wide).
for ($c = 5; $c -le 8; $c++) {
    $ws.Columns.Item($c).ColumnWidth = 7.17
}
